$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Version value from "0.1" to "1.0" (keep it stored as text,
#    matching the original cell's text type, by using a leading apostrophe)
$ws.Range("D2").Value = "'1.0"

# 2. Swap the "detalhar diária" (TC2 block) and "cancelar diária" (TC3 block)
#    content so that TC2's step/expected-result becomes the "cancelar" text
#    and TC3's step/expected-result becomes the "detalhar" text.
$stepTC2 = $ws.Range("B18").Value2
$resultTC2 = $ws.Range("D18").Value2
$stepTC3 = $ws.Range("B25").Value2
$resultTC3 = $ws.Range("D25").Value2

$ws.Range("B18").Value = $stepTC3
$ws.Range("D18").Value = $resultTC3
$ws.Range("B25").Value = $stepTC2
$ws.Range("D25").Value = $resultTC2
